# fixed some bugs in holdAndWin
# Rows 2-23 (A:F) on the active sheet get re-keyed/re-ordered to new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1001,18,30,75,60,72),
    @(301,6,45,30,60,45),
    @(601,9,60,67,60,42),
    @(1201,2,10,10,10,10),
    @(1203,3,15,15,15,15),
    @(902,1,0,0,0,0),
    @(501,9,52,30,75,45),
    @(401,9,48,67,75,45),
    @(201,9,30,15,45,30),
    @(801,3,67,65,52,45),
    @(701,3,90,45,97,15),
    @(1202,2,10,10,10,10),
    @(101,9,30,15,60,15),
    @(901,16,15,45,60,60),
    @(1101,0,15,30,30,0),
    @(802,0,4,5,4,0),
    @(1,0,2,2,2,2),
    @(502,0,4,0,0,0),
    @(3,0,3,3,3,3),
    @(2,0,2,2,2,2),
    @(602,0,0,4,0,9),
    @(402,0,0,4,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
